# Rename header labels on existing sheets
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" worksheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy the header formatting (bold, centered, bordered) from an existing sheet
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Copy the date-column number formatting down column A
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A69").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row
$wsForecast.Cells.Item(1,1).Value = "ds"
$wsForecast.Cells.Item(1,2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value = "yhat_upper"

$poForecastData = @(
    @(2, 44934.99999999999, 140, -509.268071911785, 712.25095289068),
    @(3, 44941.99999999999, 146, -491.531119671411, 777.2163807951175),
    @(4, 44955.99999999999, 158, -465.9466193307448, 785.0201682969018),
    @(5, 44962.99999999999, 164, -458.009479288362, 786.231914104674),
    @(6, 44969.99999999999, 170, -449.0512081206953, 768.835534574087),
    @(7, 44976.99999999999, 176, -497.6577358887093, 786.986147579642),
    @(8, 44983.99999999999, 182, -427.2892466958896, 794.4503047645338),
    @(9, 44990.99999999999, 188, -423.165102210582, 817.6915971324016),
    @(10, 44997.99999999999, 194, -436.2271050565666, 798.5954764017375),
    @(11, 45004.99999999999, 200, -427.594496131109, 808.1691037888921),
    @(12, 45011.99999999999, 206, -400.3411246384236, 847.402668867808),
    @(13, 45018.99999999999, 212, -391.6865867693408, 790.1340420051923),
    @(14, 45025.99999999999, 218, -409.3101970880464, 798.8647431343283),
    @(15, 45039.99999999999, 230, -382.5146302777541, 848.8764103560603),
    @(16, 45046.99999999999, 236, -360.4998310525507, 818.8410860535673),
    @(17, 45053.99999999999, 242, -397.8703710114086, 829.064493473447),
    @(18, 45060.99999999999, 248, -373.0481022033462, 839.0514536331451),
    @(19, 45067.99999999999, 254, -357.9141494888169, 866.342418004322),
    @(20, 45088.99999999999, 272, -329.2674559925852, 904.9654144095946),
    @(21, 45095.99999999999, 278, -315.9263550699239, 873.4200801761807),
    @(22, 45102.99999999999, 284, -336.1506758462365, 918.9792755930115),
    @(23, 45109.99999999999, 290, -314.6538023276768, 944.0806703975377),
    @(24, 45123.99999999999, 302, -331.6798766518637, 946.0528716741946),
    @(25, 45151.99999999999, 327, -268.964966886407, 944.1338948131227),
    @(26, 45158.99999999999, 333, -281.0871343990136, 945.6077467513425),
    @(27, 45165.99999999999, 339, -282.935663425704, 970.0563790404212),
    @(28, 45172.99999999999, 345, -301.6550056539739, 954.8717558270615),
    @(29, 45179.99999999999, 351, -229.8205563213764, 956.6508489782531),
    @(30, 45186.99999999999, 357, -239.6449178320759, 936.1099272530207),
    @(31, 45193.99999999999, 363, -228.5552322003315, 1037.81592287088),
    @(32, 45207.99999999999, 375, -268.1100296678053, 962.6374964395632),
    @(33, 45214.99999999999, 381, -241.9194848718516, 977.7000383339646),
    @(34, 45221.99999999999, 387, -249.7097474160508, 1001.758753715371),
    @(35, 45228.99999999999, 393, -217.5626093175038, 1022.57490408559),
    @(36, 45333.99999999999, 484, -138.5183560656088, 1116.991522617708),
    @(37, 45347.99999999999, 496, -136.7846062716652, 1176.418841370063),
    @(38, 45354.99999999999, 502, -107.5119296921899, 1082.278572361845),
    @(39, 45361.99999999999, 508, -149.7839561250309, 1074.808020361987),
    @(40, 45368.99999999999, 514, -111.8654180645876, 1126.066809689185),
    @(41, 45375.99999999999, 520, -75.26331783992453, 1155.560661555817),
    @(42, 45382.99999999999, 526, -98.79775179137216, 1176.448424922013),
    @(43, 45389.99999999999, 532, -58.05462879824547, 1172.974179912281),
    @(44, 45396.99999999999, 538, -64.09737373388973, 1129.469442836296),
    @(45, 45403.99999999999, 544, -71.55420556425791, 1211.400887929414),
    @(46, 45417.99999999999, 556, -96.30452366921745, 1183.603409069414),
    @(47, 45424.99999999999, 562, -72.52098184444615, 1148.408374160739),
    @(48, 45431.99999999999, 568, -86.41441264484511, 1184.573825902956),
    @(49, 45438.99999999999, 574, 11.77805006126809, 1171.437419248963),
    @(50, 45445.99999999999, 580, -32.76338718178789, 1158.358805118404),
    @(51, 45459.99999999999, 592, -1.556981587721496, 1196.599309852892),
    @(52, 45466.99999999999, 598, -30.86743696245978, 1242.403157644938),
    @(53, 45473.99999999999, 604, 18.41471679263269, 1294.925293447056),
    @(54, 45480.99999999999, 610, 28.12932776450888, 1229.982882044827),
    @(55, 45487.99999999999, 616, 18.80104238263547, 1177.130782183918),
    @(56, 45529.99999999999, 653, -24.52633651268444, 1250.443414078641),
    @(57, 45536.99999999999, 659, 46.72075478060982, 1254.759726732266),
    @(58, 45578.99999999999, 695, 91.40405642949368, 1300.613357348439),
    @(59, 45599.99999999999, 713, 41.97761620605446, 1309.863927246907),
    @(60, 45606.99999999999, 719, 116.7145101006782, 1331.684719313594),
    @(61, 45613.99999999999, 725, 86.13268544895072, 1376.476705840222),
    @(62, 45620.99999999999, 731, 109.2834753926712, 1410.03226616954),
    @(63, 45627.99999999999, 737, 107.0564663724158, 1377.02260667229),
    @(64, 45634.99999999999, 743, 127.5446912317889, 1403.374447370524),
    @(65, 45641.99999999999, 749, 141.1412178085606, 1371.112295456885),
    @(66, 45648.99999999999, 755, 163.2921733326343, 1329.315718892214),
    @(67, 45655.99999999999, 761, 141.5971799323459, 1364.964755295942),
    @(68, 45662.99999999999, 767, 146.4327947073519, 1396.194937573095),
    @(69, 45669.99999999999, 773, 170.2628451668542, 1367.334073799058)
)

foreach ($row in $poForecastData) {
    $r = $row[0]
    $wsForecast.Cells.Item($r,1).Value = $row[1]
    $wsForecast.Cells.Item($r,2).Value = $row[2]
    $wsForecast.Cells.Item($r,3).Value = $row[3]
    $wsForecast.Cells.Item($r,4).Value = $row[4]
}

Write-Host "Added PO Forecast sheet with $($poForecastData.Count) data rows"
